$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New/updated data for rows 2-13 (columns A-T), per the LR-pair recompute
# ("Natmi following Dr Hou advice"): cluster sizes E/K now reflect 3 replicate
# cells, a new "M2" target cluster is introduced, and all derived stats are
# recalculated accordingly.
$rows = @(
    @("ECs", "Mmrn2", "Cd248", "ECs", 3, 1, 41.006619, 123.019857, 0.9349081063755518, 0.9349081063755517, 3, 1, 4.076744333333333, 12.230233, 0.05973665973801366, 0.05973665973801367, 167.173501637409, 1504.561514736681, 0.05584828743686702, 0.05584828743686702),
    @("ECs", "Mmrn2", "Cd248", "FAPs", 3, 1, 41.006619, 123.019857, 0.9349081063755518, 0.9349081063755517, 3, 1, 59.98702866666667, 179.961086, 0.8789917706772566, 0.8789917706772566, 2459.865229476078, 22138.78706528471, 0.8217765318435672, 0.8217765318435671),
    @("ECs", "Mmrn2", "Cd248", "M2", 3, 1, 41.006619, 123.019857, 0.9349081063755518, 0.9349081063755517, 1, 0.3333333333333333, 0.05860766666666667, 0.175823, 0.0008587799369903073, 0.0008587799369903073, 2.403302257479, 21.629720317311, 0.0008028803246849239, 0.0008028803246849238),
    @("ECs", "Mmrn2", "Cd248", "sCs", 3, 1, 41.006619, 123.019857, 0.9349081063755518, 0.9349081063755517, 3, 1, 4.122887, 12.368661, 0.06041278964773932, 0.06041278964773932, 169.065656389053, 1521.590907501477, 0.05648040677043251, 0.0564804067704325),
    @("FAPs", "Mmrn2", "Cd248", "ECs", 3, 1, 2.068862333333333, 6.206586999999999, 0.0471679015138598, 0.04716790151385979, 3, 1, 4.076744333333333, 12.230233, 0.05973665973801366, 0.05973665973801367, 8.434222793863443, 75.90800514477098, 0.002817652883289582, 0.002817652883289582),
    @("FAPs", "Mmrn2", "Cd248", "FAPs", 3, 1, 2.068862333333333, 6.206586999999999, 0.0471679015138598, 0.04716790151385979, 3, 1, 59.98702866666667, 179.961086, 0.8789917706772566, 0.8789917706772566, 124.1049040970536, 1116.944136873482, 0.04146019727079808, 0.04146019727079807),
    @("FAPs", "Mmrn2", "Cd248", "M2", 3, 1, 2.068862333333333, 6.206586999999999, 0.0471679015138598, 0.04716790151385979, 1, 0.3333333333333333, 0.05860766666666667, 0.175823, 0.0008587799369903073, 0.0008587799369903073, 0.1212511940112222, 1.091260746101, 0.0000405068474900375404532983092398268354372703470289707183837890625, 0.0000405068474900375269007711531710214103441103361546993255615234375),
    @("FAPs", "Mmrn2", "Cd248", "sCs", 3, 1, 2.068862333333333, 6.206586999999999, 0.0471679015138598, 0.04716790151385979, 3, 1, 4.122887, 12.368661, 0.06041278964773932, 0.06041278964773932, 8.529685618889665, 76.76717057000698, 0.002849544512282097, 0.002849544512282096),
    @("sCs", "Mmrn2", "Cd248", "ECs", 3, 1, 0.7861760000000001, 2.358528, 0.01792399211058844, 0.01792399211058843, 3, 1, 4.076744333333333, 12.230233, 0.05973665973801366, 0.05973665973801367, 3.205038553002667, 28.845346977024, 0.001070719417857063, 0.001070719417857063),
    @("sCs", "Mmrn2", "Cd248", "FAPs", 3, 1, 0.7861760000000001, 2.358528, 0.01792399211058844, 0.01792399211058843, 3, 1, 59.98702866666667, 179.961086, 0.8789917706772566, 0.8789917706772566, 47.16036224904534, 424.4432602414081, 0.01575504156289131, 0.01575504156289131),
    @("sCs", "Mmrn2", "Cd248", "M2", 3, 1, 0.7861760000000001, 2.358528, 0.01792399211058844, 0.01792399211058843, 1, 0.3333333333333333, 0.05860766666666667, 0.175823, 0.0008587799369903073, 0.0008587799369903073, 0.04607594094933334, 0.4146834685440001, 0.0000153927648153458996520097368243540358889731578528881072998046875, 0.0000153927648153458996520097368243540358889731578528881072998046875),
    @("sCs", "Mmrn2", "Cd248", "sCs", 3, 1, 0.7861760000000001, 2.358528, 0.01792399211058844, 0.01792399211058843, 3, 1, 4.122887, 12.368661, 0.06041278964773932, 0.06041278964773932, 3.241314810112, 29.171833291008, 0.001082838365024718, 0.001082838365024718)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $rows[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
